$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.336.04"
$ws.Range("E2").Value = "  +1.23%  "

$ws.Range("D3").Value = "1.827.76"
$ws.Range("E3").Value = "  +0.18%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'314.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.03%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("D7").Value = "'0.4463"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.92%  "

$ws.Range("D8").Value = "'0.3758"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.75%  "

$ws.Range("D9").Value = "'0.07538"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.13%  "

$ws.Range("D10").Value = "'0.8911"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.73%  "

$ws.Range("D11").Value = "'21.02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.22%  "

$ws.Range("D12").Value = "1.832.05"
$ws.Range("E12").Value = "  +0.09%  "

$ws.Range("D13").Value = "'6.742"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.41%  "

$ws.Range("D14").Value = "'94.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.23%  "

$ws.Range("D15").Value = "'5.411"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.37%  "

$ws.Range("D16").Value = "'0.07118"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.01%  "

$ws.Range("E17").Value = "  +0.07%  "

$ws.Range("D18").Value = "'0.000008810"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.57%  "

$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("D20").Value = "'15.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.52%  "

$ws.Range("D21").Value = "27.346.30"
$ws.Range("E21").Value = "  +1.14%  "

$ws.Range("D22").Value = "'5.274"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.82%  "

$ws.Range("D23").Value = "'10.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'1.981"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'151.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("D26").Value = "'2.339"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.43%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'18.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.11%  "

$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "'5.371"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.76%  "

$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "'117.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.89%  "

$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.08833"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.37%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'0.7838"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.10%  "

$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").Value = "'1.199"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.96%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.519"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.28%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.895"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.18%  "

$ws.Range("B35").Value = "Frax"
$ws.Range("C35").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D35").Value = "'1.002"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.10%  "

$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "'1.106"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.61%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.01988"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.50%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.05323"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.98%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'7.365"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.66%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.5309"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.52%  "

$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "'2.875"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.35%  "

$ws.Range("D42").Value = "'0.1729"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.34%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'2.285"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +17.38%  "

$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "'8.742"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.23%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.5142"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.68%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'10.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.78%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'1.705"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.09%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'105.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.26%  "

$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "'1.001"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.13%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.06372"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.73%  "

$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'0.9356"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.45%  "
